$wb = $excel.ActiveWorkbook

# --- Scrum data for week6: fill in "Wat heb ik gedaan?" / "Welke problemen heb ik gehad?"
# Order matters: it determines the order new entries land in sharedStrings.xml,
# matching how the author actually typed them into the sheet.
$ws6 = $wb.Worksheets.Item("week6")

$ws6.Range("B2").Value = "mockups bijgewerkt, github in orde gebracht"
$ws6.Range("C2").Value = "kon eerst niet inloggen bij github (ondertussen wel)"
$ws6.Range("C4").Value = "security tabellen implementeren"
$ws6.Range("B3").Value = "Wcf service gemaakt, logging geïmplementeerd"
$ws6.Range("C3").Value = "/"
$ws6.Range("B4").Value = "Login en registratie begonnen"

# --- Move the active tab / selection from week5 to week6
$ws6.Activate()
$ws6.Range("C4").Select()
